$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$names = @(
    "3535 Opal Meadow Heights Aged Care Community Meadow Heights",
    "95 Napier Street Apartment Complex Fitzroy",
    "Adorn Cosmetics Clayton",
    "Al Haj Halal Meats Glenroy",
    "Al-Taqwa College Truganina",
    "Amiga Montessori Craigieburn",
    "Australia Post Distribution Centre Sunshine West",
    "Baxter Foods Australia Campbellfield",
    "Best&Less Fountain Gate Narre Warren",
    "Budget Car and Truck Rentals Campbellfield",
    "CS Square Caroline Springs",
    "Campbellfield Ford Complex Vaccination Clinic Campbellfield",
    "Cannie Road Construction Site Cannie",
    "Cannie Road Construction Site Cannie",
    "Caroline Springs Police Station",
    "Cedars Medical Clinic Coburg",
    "Chemist Warehouse Campbellfield DC",
    "Chemist Warehouse Fillo Drive Somerton",
    "City of Moreland Community",
    "City of Wyndham Community",
    "Classy Cabinets and Kitchens Craigieburn",
    "Coles Aurora Village Epping",
    "Coles Broadmeadows Central Shopping Centre",
    "Coles Campbellfield Plaza Campbellfield",
    "Coles Coburg North Village",
    "Coles Greenvale Shopping Centre",
    "Coles Pakenham Place Shopping Centre",
    "Coles Roxburgh Village Roxburgh Park",
    "Community Kids Bayswater Early Education Centre Bayswater North",
    "Community Kids Meadow Heights",
    "Construction Site Olea Apartment Caulfield North",
    "Costco Wholesale Epping",
    "Crusader Caravans Epping",
    "DayHab Rehabilitation Treatment Centre Ringwood East",
    "Direct Freight Express Cambellfield",
    "Don Watson Coldstore Derrimut",
    "Epworth Healthcare Epworth Richmond Emergency Department",
    "Fine Food Holdings Pty Ltd Dandenong South",
    "Fitzroy Community School Fitzroy North",
    "Fonterra Manufacturing Workplace Campbellfield",
    "General Foods Campbellfield",
    "Glenroy West Primary School",
    "Goodstart Early Learning Altona",
    "Green Leaves Early Learning Cairnlea",
    "Green Leaves Early Learning Centre Highlands Craigieburn",
    "Gumboots Early Learning Centre South Morang",
    "Hamilton Marino 236 Jasper Road McKinnon",
    "Health Care Providers Association South Melbourne",
    "Hello Fresh Warehouse Ravenhall",
    "IGA Meadow Heights Shopping Centre Meadow Heights",
    "ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine",
    "Ibis Kingsgate Hotel Melbourne",
    "Ilim College Glenroy Campus Hadfield",
    "Ilim Learning Sanctuary Glenroy",
    "Industrial Galvanizers Valmont Coatings Campbellfield",
    "KFC Fawkner",
    "Kasr Sweets Coolaroo",
    "Kids House Early Learning Cheltenham",
    "Kippers Seafood Werribee",
    "Kool Kidz Childcare Narre Warren",
    "Learning Nest Early Learning Centre Meadow Heights",
    "Level Crossing Removal Project Lilydale Construction Site John Street",
    "Lineage Logistics Laverton North",
    "Linfox Somerton National Distribution Centre Somerton",
    "McDonald's Craigieburn North",
    "Mecca D.C Warehouse Melbourne Airport",
    "Melbourne Assessment Prison West Melbourne",
    "Melbourne Metropolitan Remand Centre Ravenhall",
    "Melbourne Truck Repairs Campbellfield",
    "Melbourne West Police Station Docklands",
    "Mercy Hospital for Women Heidelberg",
    "Mernda YMCA Early Learning Centre Mernda",
    "Mill Park Police Station Mill Park",
    "MyCentre Childcare Broadmeadows",
    "National Gallery of Victoria Melbourne",
    "Nido Early School Ascot Vale",
    "Nido Early School Glenroy",
    "Nido Early School Moonee Ponds",
    "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
    "Northern Health The Northern Hospital Epping",
    "OnQ Plumbing and Excavations Craigieburn",
    "Oporto Coolaroo",
    "Oscar Romero Catholic Primary School Craigieburn",
    "Our Lady Help of Christian's Primary School Brunswick East",
    "Paisley Park Early Learning Centre Bundoora",
    "Panorama Construction Site Whitehorse Rd Box Hill",
    "Private Residence Northern Community Services Fawkner",
    "Ramsay Health Care Warringal Private Hospital Heidelberg",
    "Richmond Quarter 261-271 Bridge Road Construction Site Richmond",
    "Sacca's Fruit World Broadmeadows Central Shopping Centre",
    "Salta Drive Construction Site Rangedale Drainage Altona North",
    "St Margaret's Primary School OSHC Maribyrnong",
    "St Vincents Hospital Emergency Department Melbourne",
    "Sultan Halal Meats & Poultry Campbellfield",
    "Tek Foods Somerton",
    "The Homestead Child and Family Centre Roxburgh Park",
    "The Huntly-Goornong Rail Works",
    "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B",
    "The Royal Melbourne Hospital AMU Ward Parkville",
    "ThorwestenCabinets Pakenham",
    "Total Window Concepts Hoppers Crossing",
    "Unilodge College Square Student Accommodation 570 Lygon Street Carlton",
    "Wallaby Childcare Wollert",
    "Werribee Mercy Hospital Emergency Department",
    "Western Health Footscray Hospital Emergency Department",
    "Western Health Sunshine Hospital Emergency Department",
    "Woodlands Long Day Care and Kindergarten Roxburgh Park",
    "Woolworths Greenvale Lakes Roxburgh Park",
    "Yara Childcare Centre Truganina"
)

$values = @(
    27,
    5,
    5,
    27,
    9,
    26,
    5,
    5,
    5,
    6,
    11,
    14,
    5,
    6,
    11,
    34,
    6,
    11,
    5,
    6,
    9,
    6,
    7,
    12,
    31,
    6,
    13,
    9,
    11,
    11,
    17,
    31,
    23,
    5,
    13,
    5,
    6,
    8,
    33,
    10,
    13,
    6,
    11,
    5,
    11,
    5,
    13,
    7,
    5,
    6,
    11,
    6,
    10,
    9,
    22,
    6,
    7,
    9,
    6,
    6,
    5,
    7,
    8,
    9,
    5,
    9,
    5,
    11,
    7,
    7,
    5,
    5,
    7,
    17,
    9,
    6,
    11,
    13,
    55,
    8,
    19,
    11,
    6,
    10,
    6,
    13,
    5,
    9,
    11,
    6,
    7,
    11,
    9,
    5,
    20,
    10,
    5,
    18,
    27,
    14,
    6,
    12,
    11,
    11,
    6,
    9,
    5,
    5,
    8
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $names[$i]
    $ws.Cells.Item($r, 2).Value2 = $values[$i]
}

Write-Output "Updated $($names.Count) rows"